$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 395
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 95
